$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-Text "D2" "63.001.55"
Set-Text "E2" "  -0.85%  "

# Row 3 - Ethereum
Set-Text "D3" "3.048.92"
Set-Text "E3" "  -1.27%  "

# Row 4 - TetherUSD
Set-Text "E4" "  -0.14%  "

# Row 5 - BNB
Set-Text "D5" "582.42"
Set-Text "E5" "  -1.65%  "

# Row 6 - Solana
Set-Text "D6" "150.95"
Set-Text "E6" "  -2.61%  "

# Row 7 - USDC
Set-Text "E7" "  -0.04%  "

# Row 8 - XRP
Set-Text "D8" "0.533"
Set-Text "E8" "  -2.38%  "

# Row 9 - LidoStakedEther
Set-Text "D9" "3.050.69"
Set-Text "E9" "  -0.92%  "

# Row 10 - Dogecoin
Set-Text "D10" "0.152"
Set-Text "E10" "  -3.23%  "

# Row 11 - Toncoin
Set-Text "D11" "5.81"
Set-Text "E11" "  -0.60%  "

# Row 12 - Cardano
Set-Text "E12" "  -2.50%  "

# Row 13 - ShibaInu
Set-Text "D13" "0.0000233"
Set-Text "E13" "  -3.15%  "

# Row 14 - Avalanche
Set-Text "D14" "35.96"
Set-Text "E14" "  -4.16%  "

# Row 15 - TRON
Set-Text "E15" "  +1.96%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-Text "D16" "3.557.36"
Set-Text "E16" "  -1.18%  "

# Row 17 - Polkadot
Set-Text "D17" "7.12"
Set-Text "E17" "  -1.05%  "

# Row 18 - WrappedBTC
Set-Text "D18" "62.967.02"
Set-Text "E18" "  -0.87%  "

# Row 19 - WrappedEther
Set-Text "D19" "3.051.15"
Set-Text "E19" "  -1.07%  "

# Row 20 - BitcoinCash
Set-Text "D20" "479.65"
Set-Text "E20" "  +0.63%  "

# Row 21 - Chainlink
Set-Text "D21" "14.28"
Set-Text "E21" "  -2.69%  "

# Row 22 - Polygon
Set-Text "D22" "0.705"
Set-Text "E22" "  -1.87%  "

# Row 23 - Uniswap
Set-Text "D23" "7.50"
Set-Text "E23" "  -0.97%  "

# Row 24 - Fetch.AI
Set-Text "E24" "  -1.15%  "

# Row 25 - Litecoin
Set-Text "D25" "81.83"
Set-Text "E25" "  +0.69%  "

# Row 26 - InternetComputer(DFINITY)
Set-Text "D26" "12.62"
Set-Text "E26" "  -2.31%  "

# Row 27 - RenderToken
Set-Text "D27" "10.55"

# Row 28 - Dai
Set-Text "E28" "  +0.00%  "

# Row 29 - NEARProtocol
Set-Text "E29" "  -0.33%  "

# Row 30 - FirstDigitalUSD
Set-Text "E30" "  -0.14%  "

# Row 31 - PancakeSwap
Set-Text "E31" "  -1.59%  "

# Row 32 - ImmutableX
Set-Text "D32" "2.19"
Set-Text "E32" "  +0.21%  "

# Row 33 - EthereumClassic
Set-Text "E33" "  +1.65%  "

# Row 34 - Hedera
Set-Text "E34" "  -3.37%  "

# Row 35 - Mantle
Set-Text "E35" "  +0.56%  "

# Row 36 - PEPE
Set-Text "E36" "  -4.62%  "

# Row 37 - Filecoin
Set-Text "E37" "  -3.37%  "

# Row 38 - Stacks
Set-Text "E38" "  -1.91%  "

# Row 39 - dogwifhat
Set-Text "D39" "3.13"
Set-Text "E39" "  -7.44%  "

# Row 40 - was Cosmos, now OKB
Set-Text "B40" "OKB"
Set-Text "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-Text "D40" "50.33"
Set-Text "E40" "  -0.96%  "

# Row 41 - was OKB, now Cosmos
Set-Text "B41" "Cosmos"
Set-Text "C41" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Text "D41" "9.16"
Set-Text "E41" "  -2.24%  "

# Row 42 - Bittensor
Set-Text "D42" "425.42"
Set-Text "E42" "  -4.32%  "

# Row 43 - was TheGraph, now Kaspa
Set-Text "B43" "Kaspa"
Set-Text "C43" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-Text "D43" "0.115"
Set-Text "E43" "  +3.49%  "

# Row 44 - was Kaspa, now TheGraph
Set-Text "B44" "TheGraph"
Set-Text "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-Text "D44" "0.285"
Set-Text "E44" "  +0.09%  "

# Row 45 - Maker
Set-Text "D45" "2.844.96"
Set-Text "E45" "  +1.53%  "

# Row 46 - VeChain
Set-Text "D46" "0.0360"
Set-Text "E46" "  -0.71%  "

# Row 47 - Arweave
Set-Text "D47" "37.74"
Set-Text "E47" "  -5.82%  "

# Row 48 - Monero
Set-Text "D48" "126.83"
Set-Text "E48" "  -3.58%  "

# Row 49 - USDe
Set-Text "D49" "1.00"

# Row 50 - InjectiveProtocol
Set-Text "D50" "24.97"
Set-Text "E50" "  -2.50%  "

# Row 51 - Stellar
Set-Text "E51" "  -1.09%  "
